$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data: id=8, name="arctique"
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "arctique"

# Scroll the view so that row 4 is the top-left visible row, matching the diff
$ws.Application.ActiveWindow.ScrollRow = 4

# Keep selection on B9 as in the original/target view
$ws.Range("B9").Select()
